# "Penalty Reward System" (unfinished) edit:
#  - Weekly Quantity: drop the week of 2024-03-18 (old row 14), which
#    shifts every following row up by one, and adjust the quantity for
#    the week of 2024-03-11 (now row 13) from 160 to 130.
#  - Monthly Trend: adjust the March 2024 (April row, A6) quantity
#    from 415 to 355.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")

# Deleting the entire row shifts rows 15:34 up into 14:33, which is what
# turns the old A1:B34 range into A1:B33.
$wsWeekly.Rows.Item(14).Delete()

# Former row 13 keeps its date (2024-03-11) but its quantity changes.
$wsWeekly.Range("B13").Value = 130

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B6").Value = 355
